$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 16,20
$data[0,0] = "ECs"
$data[0,1] = "Gpc3"
$data[0,2] = "Cd81"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 1.465847333333333
$data[0,7] = 4.397542
$data[0,8] = 0.0219775736133859
$data[0,9] = 0.02197757361338591
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 181.5963696666667
$data[0,13] = 544.7891090000001
$data[0,14] = 0.361289261089317
$data[0,15] = 0.361289261089317
$data[0,16] = 266.1925542188976
$data[0,17] = 2395.732987970078
$data[0,18] = 0.007940261331316265
$data[0,19] = 0.007940261331316267
$data[1,0] = "ECs"
$data[1,1] = "Gpc3"
$data[1,2] = "Cd81"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 1.465847333333333
$data[1,7] = 4.397542
$data[1,8] = 0.0219775736133859
$data[1,9] = 0.02197757361338591
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 203.8031973333333
$data[1,13] = 611.409592
$data[1,14] = 0.4054701462774667
$data[1,15] = 0.4054701462774667
$data[1,16] = 298.7443733358737
$data[1,17] = 2688.699360022864
$data[1,18] = 0.008911249987843374
$data[1,19] = 0.008911249987843376
$data[2,0] = "ECs"
$data[2,1] = "Gpc3"
$data[2,2] = "Cd81"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 1.465847333333333
$data[2,7] = 4.397542
$data[2,8] = 0.0219775736133859
$data[2,9] = 0.02197757361338591
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 28.96574966666667
$data[2,13] = 86.89724900000002
$data[2,14] = 0.0576278827224213
$data[2,15] = 0.0576278827224213
$data[2,16] = 42.45936690688423
$data[2,17] = 382.134302161958
$data[2,18] = 0.001266521034715584
$data[2,19] = 0.001266521034715584
$data[3,0] = "ECs"
$data[3,1] = "Gpc3"
$data[3,2] = "Cd81"
$data[3,3] = "sCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 1.465847333333333
$data[3,7] = 4.397542
$data[3,8] = 0.0219775736133859
$data[3,9] = 0.02197757361338591
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 88.26896899999998
$data[3,13] = 264.806907
$data[3,14] = 0.175612709910795
$data[3,15] = 0.175612709910795
$data[3,16] = 129.3888328247326
$data[3,17] = 1164.499495422594
$data[3,18] = 0.003859541259510681
$data[3,19] = 0.003859541259510682
$data[4,0] = "FAPs"
$data[4,1] = "Gpc3"
$data[4,2] = "Cd81"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 51.92481233333333
$data[4,7] = 155.774437
$data[4,8] = 0.7785131230699432
$data[4,9] = 0.7785131230699435
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 181.5963696666667
$data[4,13] = 544.7891090000001
$data[4,14] = 0.361289261089317
$data[4,15] = 0.361289261089317
$data[4,16] = 9429.357415356291
$data[4,17] = 84864.21673820663
$data[4,18] = 0.2812684309822763
$data[4,19] = 0.2812684309822764
$data[5,0] = "FAPs"
$data[5,1] = "Gpc3"
$data[5,2] = "Cd81"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 51.92481233333333
$data[5,7] = 155.774437
$data[5,8] = 0.7785131230699432
$data[5,9] = 0.7785131230699435
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 203.8031973333333
$data[5,13] = 611.409592
$data[5,14] = 0.4054701462774667
$data[5,15] = 0.4054701462774667
$data[5,16] = 10582.44277446663
$data[5,17] = 95241.98497019969
$data[5,18] = 0.3156638298900973
$data[5,19] = 0.3156638298900974
$data[6,0] = "FAPs"
$data[6,1] = "Gpc3"
$data[6,2] = "Cd81"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 51.92481233333333
$data[6,7] = 155.774437
$data[6,8] = 0.7785131230699432
$data[6,9] = 0.7785131230699435
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 28.96574966666667
$data[6,13] = 86.89724900000002
$data[6,14] = 0.0576278827224213
$data[6,15] = 0.0576278827224213
$data[6,16] = 1504.041115535979
$data[6,17] = 13536.37003982381
$data[6,18] = 0.04486406295414063
$data[6,19] = 0.04486406295414064
$data[7,0] = "FAPs"
$data[7,1] = "Gpc3"
$data[7,2] = "Cd81"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 51.92481233333333
$data[7,7] = 155.774437
$data[7,8] = 0.7785131230699432
$data[7,9] = 0.7785131230699435
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 88.26896899999998
$data[7,13] = 264.806907
$data[7,14] = 0.175612709910795
$data[7,15] = 0.175612709910795
$data[7,16] = 4583.349650181816
$data[7,17] = 41250.14685163635
$data[7,18] = 0.136716799243429
$data[7,19] = 0.136716799243429
$data[8,0] = "M2"
$data[8,1] = "Gpc3"
$data[8,2] = "Cd81"
$data[8,3] = "ECs"
$data[8,4] = 2
$data[8,5] = 0.6666666666666666
$data[8,6] = 0.05191500000000001
$data[8,7] = 0.155745
$data[8,8] = 0.0007783660059225787
$data[8,9] = 0.0007783660059225788
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 181.5963696666667
$data[8,13] = 544.7891090000001
$data[8,14] = 0.361289261089317
$data[8,15] = 0.361289261089317
$data[8,16] = 9.427575531245003
$data[8,17] = 84.84817978120502
$data[8,18] = 0.0002812152791368114
$data[8,19] = 0.0002812152791368115
$data[9,0] = "M2"
$data[9,1] = "Gpc3"
$data[9,2] = "Cd81"
$data[9,3] = "FAPs"
$data[9,4] = 2
$data[9,5] = 0.6666666666666666
$data[9,6] = 0.05191500000000001
$data[9,7] = 0.155745
$data[9,8] = 0.0007783660059225787
$data[9,9] = 0.0007783660059225788
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 203.8031973333333
$data[9,13] = 611.409592
$data[9,14] = 0.4054701462774667
$data[9,15] = 0.4054701462774667
$data[9,16] = 10.58044298956
$data[9,17] = 95.22398690604001
$data[9,18] = 0.0003156041782788355
$data[9,19] = 0.0003156041782788355
$data[10,0] = "M2"
$data[10,1] = "Gpc3"
$data[10,2] = "Cd81"
$data[10,3] = "M2"
$data[10,4] = 2
$data[10,5] = 0.6666666666666666
$data[10,6] = 0.05191500000000001
$data[10,7] = 0.155745
$data[10,8] = 0.0007783660059225787
$data[10,9] = 0.0007783660059225788
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 28.96574966666667
$data[10,13] = 86.89724900000002
$data[10,14] = 0.0576278827224213
$data[10,15] = 0.0576278827224213
$data[10,16] = 1.503756893945001
$data[10,17] = 13.53381204550501
$data[10,18] = [double]"4.485558490442585E-05"
$data[10,19] = [double]"4.485558490442585E-05"
$data[11,0] = "M2"
$data[11,1] = "Gpc3"
$data[11,2] = "Cd81"
$data[11,3] = "sCs"
$data[11,4] = 2
$data[11,5] = 0.6666666666666666
$data[11,6] = 0.05191500000000001
$data[11,7] = 0.155745
$data[11,8] = 0.0007783660059225787
$data[11,9] = 0.0007783660059225788
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 88.26896899999998
$data[11,13] = 264.806907
$data[11,14] = 0.175612709910795
$data[11,15] = 0.175612709910795
$data[11,16] = 4.582483525635
$data[11,17] = 41.242351730715
$data[11,18] = 0.0001366909636025059
$data[11,19] = 0.0001366909636025059
$data[12,0] = "sCs"
$data[12,1] = "Gpc3"
$data[12,2] = "Cd81"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 13.25483966666666
$data[12,7] = 39.76451899999999
$data[12,8] = 0.1987309373107482
$data[12,9] = 0.1987309373107483
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 181.5963696666667
$data[12,13] = 544.7891090000001
$data[12,14] = 0.361289261089317
$data[12,15] = 0.361289261089317
$data[12,16] = 2407.030763980396
$data[12,17] = 21663.27687582357
$data[12,18] = 0.0717993534965876
$data[12,19] = 0.07179935349658761
$data[13,0] = "sCs"
$data[13,1] = "Gpc3"
$data[13,2] = "Cd81"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 13.25483966666666
$data[13,7] = 39.76451899999999
$data[13,8] = 0.1987309373107482
$data[13,9] = 0.1987309373107483
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 203.8031973333333
$data[13,13] = 611.409592
$data[13,14] = 0.4054701462774667
$data[13,15] = 0.4054701462774667
$data[13,16] = 2701.37870420736
$data[13,17] = 24312.40833786624
$data[13,18] = 0.08057946222124714
$data[13,19] = 0.08057946222124716
$data[14,0] = "sCs"
$data[14,1] = "Gpc3"
$data[14,2] = "Cd81"
$data[14,3] = "M2"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 13.25483966666666
$data[14,7] = 39.76451899999999
$data[14,8] = 0.1987309373107482
$data[14,9] = 0.1987309373107483
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 28.96574966666667
$data[14,13] = 86.89724900000002
$data[14,14] = 0.0576278827224213
$data[14,15] = 0.0576278827224213
$data[14,16] = 383.9363676564701
$data[14,17] = 3455.427308908231
$data[14,18] = 0.01145244314866066
$data[14,19] = 0.01145244314866066
$data[15,0] = "sCs"
$data[15,1] = "Gpc3"
$data[15,2] = "Cd81"
$data[15,3] = "sCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 13.25483966666666
$data[15,7] = 39.76451899999999
$data[15,8] = 0.1987309373107482
$data[15,9] = 0.1987309373107483
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 88.26896899999998
$data[15,13] = 264.806907
$data[15,14] = 0.175612709910795
$data[15,15] = 0.175612709910795
$data[15,16] = 1169.99103163697
$data[15,17] = 10529.91928473273
$data[15,18] = 0.0348996784442528
$data[15,19] = 0.03489967844425281
$ws.Range("A2:T17").Value2 = $data
Write-Host "done"
